$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.285.40'
$ws.Range("E2").Value = '  -5.27%  '
$ws.Range("D3").Value = '3.489.42'
$ws.Range("E3").Value = '  -6.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.66%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.95%  '
$ws.Range("D7").Value = '3.479.36'
$ws.Range("E7").Value = '  -6.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.596'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.79%  '
$ws.Range("E9").Value = '  +0.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.653'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.31'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -8.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.139'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -14.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000241'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -18.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.55'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -11.03%  '
$ws.Range("D15").Value = '4.086.28'
$ws.Range("E15").Value = '  -5.16%  '
$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.125'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.07%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.509.98'
$ws.Range("E17").Value = '  -5.64%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '65.357.39'
$ws.Range("E18").Value = '  -4.95%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -8.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.81'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -9.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.03'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -9.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '383.46'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -11.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.80'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.69%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.98'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.62%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.99'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -9.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.49'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -10.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.67'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -10.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.20'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.55'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -11.07%  '
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '605.01'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.46%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '64.44'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.29%  '
$ws.Range("B35").Value = 'Cosmos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.69'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.109'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -9.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '40.26'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -8.81%  '
$ws.Range("E38").Value = '  +0.12%  '
$ws.Range("E39").Value = '  +0.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.362'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -11.58%  '
$ws.Range("D41").Value = '0.0₃0715'
$ws.Range("E41").Value = '  -19.93%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.127'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.98%  '
$ws.Range("D43").Value = '2.819.60'
$ws.Range("E43").Value = '  +0.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.72'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -11.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0397'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -10.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -12.83%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.127'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.40%  '
$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '135.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.41'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -11.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.05'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -13.07%  '
